$wb = $excel.ActiveWorkbook

# --- Repayment schedule: add O3:O14 = 0 and an (empty, styled) P2 cell ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$p2 = $wsRepay.Range("P2")
$p2.VerticalAlignment = -4108   # xlCenter
$p2.WrapText = $true

for ($r = 3; $r -le 14; $r++) {
    $cell = $wsRepay.Range("O$r")
    $cell.Value = 0
    $cell.VerticalAlignment = -4108   # xlCenter
    $cell.WrapText = $true
}

# --- Update selections on each sheet (order matters: selecting a range
#     activates that sheet's tab, so the sheet that should end up active
#     must be activated/selected last) ---

$wsLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsLoanInput.Range("B21").Select()

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A7:XFD14").Select()

$wsRepay.Range("A15:XFD15").Select()

$wsFloating = $wb.Worksheets.Item("Floating Interest Rates")
$wsFloating.Range("A3:B6").Select()

$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("C2").Select()

# Transactions becomes the active tab/sheet.
$wsTransactions.Activate()
